$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct text assignments (values that are unambiguously text, e.g.
# contain multiple "." separators or percent signs, so Excel keeps them as text)
$ws.Range("D2").Value = '42.427.58'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '2.241.00'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("E6").Value = '  -1.47%  '
$ws.Range("E7").Value = '  -1.89%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -2.21%  '
$ws.Range("E10").Value = '  -3.13%  '
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("E12").Value = '  -4.29%  '
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("D14").Value = '2.576.19'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("E16").Value = '  -2.48%  '
$ws.Range("D17").Value = '2.218.95'
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("D18").Value = '42.264.62'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("E19").Value = '  +5.15%  '
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("E22").Value = '  +2.91%  '
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("E24").Value = '  -6.54%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  -2.58%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("E29").Value = '  -2.21%  '
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("E34").Value = '  -3.38%  '
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("E36").Value = '  -8.74%  '
$ws.Range("E37").Value = '  -8.42%  '
$ws.Range("E38").Value = '  -3.60%  '
$ws.Range("E39").Value = '  -5.74%  '
$ws.Range("E40").Value = '  -2.13%  '
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("E43").Value = '  -1.86%  '
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("E45").Value = '  -2.72%  '
$ws.Range("E46").Value = '  -1.93%  '
$ws.Range("E47").Value = '  -1.85%  '
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("E49").Value = '  -1.63%  '
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").Value = '2.447.72'
$ws.Range("E51").Value = '  -0.69%  '

# Values that look like plain numbers (e.g. "0.622") must be forced to stay as
# TEXT (matching the source workbook, which stores them as inline strings) without
# altering the target cell's style/number-format. We stage the text in a scratch
# cell formatted as Text ("@"), copy it, and paste-special values-only into the
# destination - PasteSpecial(xlPasteValues) carries over the TEXT cell type from
# the source without touching the destination's own style/number-format.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = '243.86'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '0.622'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '74.59'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '43.05'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '0.0964'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '14.39'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '0.843'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '0.0000107'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '73.23'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '11.32'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '231.71'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '11.50'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '167.17'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '20.67'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '5.72'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '30.70'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '0.110'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '4.40'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '0.0306'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '13.41'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '5.71'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '65.08'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '8.77'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Value = '105.27'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").Clear() | Out-Null
